$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815299
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45144.625
$arr[0,3] = 'FK Decic Tuzi'
$arr[0,4] = 'OFK Petrovac'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'A'
$arr[0,10] = 1.95
$arr[0,11] = 3.3
$arr[0,12] = 3.4
$arr[0,13] = 1.909
$arr[0,14] = 3.3
$arr[0,15] = 3.5
$arr[0,16] = -0.5
$arr[0,17] = 1.975
$arr[0,18] = 1.825
$arr[0,19] = 2.25
$arr[0,20] = 1.925
$arr[0,21] = 1.875
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 2.5
$arr[0,25] = -1
$arr[0,26] = 0.825
$arr[0,27] = 0.925
$arr[0,28] = -1
$ws.Range("B14:AD14").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815302
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45144.625
$arr[0,3] = 'Sutjeska Niksic'
$arr[0,4] = 'FK Arsenal'
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 2
$arr[0,9] = 'D'
$arr[0,10] = 1.533
$arr[0,11] = 3.6
$arr[0,12] = 5.75
$arr[0,13] = 1.533
$arr[0,14] = 3.6
$arr[0,15] = 5.75
$arr[0,16] = -1
$arr[0,17] = 1.95
$arr[0,18] = 1.85
$arr[0,19] = 2.25
$arr[0,20] = 1.925
$arr[0,21] = 1.875
$arr[0,22] = -1
$arr[0,23] = 2.6
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.8500000000000001
$arr[0,27] = 0.925
$arr[0,28] = -1
$ws.Range("B15:AD15").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815331
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45196.54166666666
$arr[0,3] = 'FK Jedinstvo Bijelo Polje'
$arr[0,4] = 'FK Decic Tuzi'
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'A'
$arr[0,10] = 3
$arr[0,11] = 3.1
$arr[0,12] = 2.2
$arr[0,13] = 5.75
$arr[0,14] = 3.6
$arr[0,15] = 1.5
$arr[0,16] = 1
$arr[0,17] = 1.825
$arr[0,18] = 1.975
$arr[0,19] = 2.25
$arr[0,20] = 1.95
$arr[0,21] = 1.85
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.5
$arr[0,25] = 0
$arr[0,26] = 0
$arr[0,27] = -1
$arr[0,28] = 0.8500000000000001
$ws.Range("B48:AD48").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815333
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45196.54166666666
$arr[0,3] = 'Sutjeska Niksic'
$arr[0,4] = 'FK Jezero'
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'D'
$arr[0,10] = 1.5
$arr[0,11] = 3.6
$arr[0,12] = 6
$arr[0,13] = 1.5
$arr[0,14] = 3.6
$arr[0,15] = 5.75
$arr[0,16] = -1
$arr[0,17] = 1.9
$arr[0,18] = 1.9
$arr[0,19] = 2.25
$arr[0,20] = 1.85
$arr[0,21] = 1.95
$arr[0,22] = -1
$arr[0,23] = 2.6
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$arr[0,27] = -0.5
$arr[0,28] = 0.475
$ws.Range("B49:AD49").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815338
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45206.58333333334
$arr[0,3] = 'OFK Petrovac'
$arr[0,4] = 'FK Decic Tuzi'
$arr[0,5] = 2
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 2
$arr[0,9] = 'A'
$arr[0,10] = 2.625
$arr[0,11] = 2.875
$arr[0,12] = 2.6
$arr[0,13] = 3.1
$arr[0,14] = 2.9
$arr[0,15] = 2.25
$arr[0,16] = 0.25
$arr[0,17] = 1.8
$arr[0,18] = 2
$arr[0,19] = 2.25
$arr[0,20] = 1.975
$arr[0,21] = 1.725
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 1.25
$arr[0,25] = -1
$arr[0,26] = 1
$arr[0,27] = 0.9750000000000001
$arr[0,28] = -1
$ws.Range("B59:AD59").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815427
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45206.58333333334
$arr[0,3] = 'FK Mornar Bar'
$arr[0,4] = 'OFK Mladost DG'
$arr[0,5] = 2
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = 'H'
$arr[0,10] = 1.833
$arr[0,11] = 3.1
$arr[0,12] = 4
$arr[0,13] = 1.833
$arr[0,14] = 3.1
$arr[0,15] = 4
$arr[0,16] = -0.5
$arr[0,17] = 1.875
$arr[0,18] = 1.925
$arr[0,19] = 2
$arr[0,20] = 1.775
$arr[0,21] = 2.025
$arr[0,22] = 0.833
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.875
$arr[0,26] = -1
$arr[0,27] = 0.7749999999999999
$arr[0,28] = -1
$ws.Range("B60:AD60").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815378
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45346.41666666666
$arr[0,3] = 'FK Rudar Pljevlja'
$arr[0,4] = 'Buducnost Podgorica'
$arr[0,5] = 3
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 0
$arr[0,9] = 'H'
$arr[0,10] = 10
$arr[0,11] = 5.5
$arr[0,12] = 1.2
$arr[0,13] = 11
$arr[0,14] = 4.75
$arr[0,15] = 1.222
$arr[0,16] = 1.75
$arr[0,17] = 1.85
$arr[0,18] = 1.95
$arr[0,19] = 2.5
$arr[0,20] = 1.825
$arr[0,21] = 1.975
$arr[0,22] = 10
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.8500000000000001
$arr[0,26] = -1
$arr[0,27] = 0.825
$arr[0,28] = -1
$ws.Range("B100:AD100").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815433
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45346.41666666666
$arr[0,3] = 'OFK Mladost DG'
$arr[0,4] = 'FK Mornar Bar'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'A'
$arr[0,10] = 2.5
$arr[0,11] = 3
$arr[0,12] = 2.6
$arr[0,13] = 2.5
$arr[0,14] = 3
$arr[0,15] = 2.6
$arr[0,16] = 0
$arr[0,17] = 1.85
$arr[0,18] = 1.95
$arr[0,19] = 2
$arr[0,20] = 1.975
$arr[0,21] = 1.825
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 1.6
$arr[0,25] = -1
$arr[0,26] = 0.95
$arr[0,27] = 0.9750000000000001
$arr[0,28] = -1
$ws.Range("B101:AD101").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815382
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45350.45833333334
$arr[0,3] = 'FK Jedinstvo Bijelo Polje'
$arr[0,4] = 'Sutjeska Niksic'
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = 0
$arr[0,8] = 1
$arr[0,9] = 'A'
$arr[0,10] = 3.75
$arr[0,11] = 3.3
$arr[0,12] = 1.833
$arr[0,13] = 4.75
$arr[0,14] = 3
$arr[0,15] = 1.75
$arr[0,16] = 0.5
$arr[0,17] = 1.975
$arr[0,18] = 1.825
$arr[0,19] = 2
$arr[0,20] = 1.85
$arr[0,21] = 1.95
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.75
$arr[0,25] = -1
$arr[0,26] = 0.825
$arr[0,27] = -1
$arr[0,28] = 0.95
$ws.Range("B105:AD105").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815434
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45350.45833333334
$arr[0,3] = 'OFK Mladost DG'
$arr[0,4] = 'FK Decic Tuzi'
$arr[0,5] = 0
$arr[0,6] = 3
$arr[0,7] = 0
$arr[0,8] = 1
$arr[0,9] = 'A'
$arr[0,10] = 4.6
$arr[0,11] = 3.6
$arr[0,12] = 1.615
$arr[0,13] = 8
$arr[0,14] = 4.75
$arr[0,15] = 1.25
$arr[0,16] = 0.75
$arr[0,17] = 1.925
$arr[0,18] = 1.875
$arr[0,19] = 2.25
$arr[0,20] = 1.9
$arr[0,21] = 1.9
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.25
$arr[0,25] = -1
$arr[0,26] = 0.875
$arr[0,27] = 0.8999999999999999
$arr[0,28] = -1
$ws.Range("B106:AD106").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815404
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45381.45833333334
$arr[0,3] = 'Sutjeska Niksic'
$arr[0,4] = 'FK Rudar Pljevlja'
$arr[0,5] = 2
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = 'H'
$arr[0,10] = 1.5
$arr[0,11] = 3.75
$arr[0,12] = 5.75
$arr[0,13] = 1.5
$arr[0,14] = 3.75
$arr[0,15] = 5.75
$arr[0,16] = -1
$arr[0,17] = 1.85
$arr[0,18] = 1.95
$arr[0,19] = 2.25
$arr[0,20] = 1.9
$arr[0,21] = 1.9
$arr[0,22] = 0.5
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0
$arr[0,26] = 0
$arr[0,27] = 0.8999999999999999
$arr[0,28] = -1
$ws.Range("B130:AD130").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815405
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45381.45833333334
$arr[0,3] = 'FK Arsenal'
$arr[0,4] = 'FK Decic Tuzi'
$arr[0,5] = 0
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'A'
$arr[0,10] = 4
$arr[0,11] = 3.2
$arr[0,12] = 1.8
$arr[0,13] = 5
$arr[0,14] = 3.25
$arr[0,15] = 1.65
$arr[0,16] = 0.75
$arr[0,17] = 1.875
$arr[0,18] = 1.925
$arr[0,19] = 2
$arr[0,20] = 1.875
$arr[0,21] = 1.925
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.6499999999999999
$arr[0,25] = -1
$arr[0,26] = 0.925
$arr[0,27] = 0
$arr[0,28] = 0
$ws.Range("B131:AD131").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815406
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45381.45833333334
$arr[0,3] = 'FK Jedinstvo Bijelo Polje'
$arr[0,4] = 'OFK Mladost DG'
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 'D'
$arr[0,10] = 2.25
$arr[0,11] = 3.3
$arr[0,12] = 2.7
$arr[0,13] = 2.05
$arr[0,14] = 3.4
$arr[0,15] = 3
$arr[0,16] = -0.25
$arr[0,17] = 1.825
$arr[0,18] = 1.975
$arr[0,19] = 2.5
$arr[0,20] = 2
$arr[0,21] = 1.8
$arr[0,22] = -1
$arr[0,23] = 2.4
$arr[0,24] = -1
$arr[0,25] = -0.5
$arr[0,26] = 0.4875
$arr[0,27] = -1
$arr[0,28] = 0.8
$ws.Range("B132:AD132").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6815403
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45381.45833333334
$arr[0,3] = 'FK Jezero'
$arr[0,4] = 'FK Mornar Bar'
$arr[0,5] = 3
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 0
$arr[0,9] = 'H'
$arr[0,10] = 2.875
$arr[0,11] = 2.75
$arr[0,12] = 2.5
$arr[0,13] = 3.1
$arr[0,14] = 2.55
$arr[0,15] = 2.55
$arr[0,16] = 0
$arr[0,17] = 2.05
$arr[0,18] = 1.75
$arr[0,19] = 1.75
$arr[0,20] = 1.95
$arr[0,21] = 1.85
$arr[0,22] = 2.1
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 1.05
$arr[0,26] = -1
$arr[0,27] = 0.95
$arr[0,28] = -1
$ws.Range("B133:AD133").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 8043518
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45385.41666666666
$arr[0,3] = 'FK Arsenal'
$arr[0,4] = 'FK Rudar Pljevlja'
$arr[0,5] = 4
$arr[0,6] = 2
$arr[0,7] = 3
$arr[0,8] = 1
$arr[0,9] = 'H'
$arr[0,10] = 1.909
$arr[0,11] = 3
$arr[0,12] = 3.9
$arr[0,13] = 1.65
$arr[0,14] = 3.3
$arr[0,15] = 5
$arr[0,16] = -0.75
$arr[0,17] = 1.875
$arr[0,18] = 1.925
$arr[0,19] = 2
$arr[0,20] = 1.8
$arr[0,21] = 2
$arr[0,22] = 0.6499999999999999
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.875
$arr[0,26] = -1
$arr[0,27] = 0.8
$arr[0,28] = -1
$ws.Range("B135:AD135").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 8043517
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45385.41666666666
$arr[0,3] = 'FK Jedinstvo Bijelo Polje'
$arr[0,4] = 'FK Decic Tuzi'
$arr[0,5] = 0
$arr[0,6] = 2
$arr[0,7] = 0
$arr[0,8] = 1
$arr[0,9] = 'A'
$arr[0,10] = 5.5
$arr[0,11] = 3.2
$arr[0,12] = 1.615
$arr[0,13] = 6.5
$arr[0,14] = 3.4
$arr[0,15] = 1.533
$arr[0,16] = 1
$arr[0,17] = 1.85
$arr[0,18] = 1.95
$arr[0,19] = 2.25
$arr[0,20] = 2.025
$arr[0,21] = 1.775
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.5329999999999999
$arr[0,25] = -1
$arr[0,26] = 0.95
$arr[0,27] = -0.5
$arr[0,28] = 0.3875
$ws.Range("B136:AD136").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6816714
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45416.5
$arr[0,3] = 'FK Arsenal'
$arr[0,4] = 'Buducnost Podgorica'
$arr[0,5] = 0
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 1
$arr[0,9] = 'A'
$arr[0,10] = 5.75
$arr[0,11] = 4.2
$arr[0,12] = 1.4
$arr[0,13] = 5.75
$arr[0,14] = 4.2
$arr[0,15] = 1.4
$arr[0,16] = 1.25
$arr[0,17] = 1.85
$arr[0,18] = 1.95
$arr[0,19] = 2.75
$arr[0,20] = 1.95
$arr[0,21] = 1.85
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.3999999999999999
$arr[0,25] = -1
$arr[0,26] = 0.95
$arr[0,27] = 0.95
$arr[0,28] = -1
$ws.Range("B158:AD158").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6816713
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45416.5
$arr[0,3] = 'FK Jedinstvo Bijelo Polje'
$arr[0,4] = 'FK Mornar Bar'
$arr[0,5] = 3
$arr[0,6] = 4
$arr[0,7] = 2
$arr[0,8] = 3
$arr[0,9] = 'A'
$arr[0,10] = 3.25
$arr[0,11] = 2.875
$arr[0,12] = 2.2
$arr[0,13] = 3
$arr[0,14] = 3
$arr[0,15] = 2.3
$arr[0,16] = 0.25
$arr[0,17] = 1.775
$arr[0,18] = 2.025
$arr[0,19] = 2
$arr[0,20] = 1.825
$arr[0,21] = 1.975
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 1.3
$arr[0,25] = -1
$arr[0,26] = 1.025
$arr[0,27] = 0.825
$arr[0,28] = -1
$ws.Range("B159:AD159").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6817582
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45424.58333333334
$arr[0,3] = 'OFK Mladost DG'
$arr[0,4] = 'Sutjeska Niksic'
$arr[0,5] = 3
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 0
$arr[0,9] = 'H'
$arr[0,10] = 4.2
$arr[0,11] = 3.3
$arr[0,12] = 1.75
$arr[0,13] = 4.2
$arr[0,14] = 3.3
$arr[0,15] = 1.75
$arr[0,16] = 0.5
$arr[0,17] = 2
$arr[0,18] = 1.8
$arr[0,19] = 2.25
$arr[0,20] = 1.8
$arr[0,21] = 2
$arr[0,22] = 3.2
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 1
$arr[0,26] = -1
$arr[0,27] = 0.8
$arr[0,28] = -1
$ws.Range("B160:AD160").Value = $arr

$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 6817581
$arr[0,1] = 'Montenegro Prva Liga'
$arr[0,2] = 45424.58333333334
$arr[0,3] = 'FK Decic Tuzi'
$arr[0,4] = 'FK Jezero'
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = 1
$arr[0,8] = 0
$arr[0,9] = 'H'
$arr[0,10] = 1.6
$arr[0,11] = 3.5
$arr[0,12] = 5
$arr[0,13] = 1.285
$arr[0,14] = 4.75
$arr[0,15] = 9.5
$arr[0,16] = -1.5
$arr[0,17] = 1.925
$arr[0,18] = 1.875
$arr[0,19] = 2.25
$arr[0,20] = 1.925
$arr[0,21] = 1.875
$arr[0,22] = 0.2849999999999999
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.925
$arr[0,26] = -1
$arr[0,27] = -0.5
$arr[0,28] = 0.4375
$ws.Range("B162:AD162").Value = $arr

